$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# A brand-new September entry was logged ("axis" at 2024-09-07 08:31:28).
# This pushes all the existing rows from row 35 downward by one row
# (dimension grows from Y88 to Y89).
$ws.Rows("35").Insert()

$ws.Range("R35").Value = "axis"
$ws.Range("S35").Value = "2024-09-07 08:31:28"
